$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.246.78'
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").Value = '1.858.08'
$ws.Range("E3").Value = '  -2.26%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4791'
$ws.Range("E7").Value = '  -2.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2800'
$ws.Range("E8").Value = '  -4.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06464'
$ws.Range("E9").Value = '  -3.22%  '

$ws.Range("D10").Value = '1.859.27'
$ws.Range("E10").Value = '  -2.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07381'
$ws.Range("E11").Value = '  +0.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.23'
$ws.Range("E12").Value = '  -4.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.018'
$ws.Range("E13").Value = '  -3.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.13'
$ws.Range("E14").Value = '  -1.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6453'
$ws.Range("E15").Value = '  -3.29%  '

$ws.Range("D16").Value = '30.188.26'
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.13'
$ws.Range("E18").Value = '  -2.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007562'
$ws.Range("E19").Value = '  -3.98%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.100.27'
$ws.Range("E20").Value = '  -2.12%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '222.20'
$ws.Range("E21").Value = '  +14.04%  '

$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.280'
$ws.Range("E23").Value = '  -2.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.074'
$ws.Range("E24").Value = '  -0.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.197'
$ws.Range("E25").Value = '  -3.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.80'
$ws.Range("E26").Value = '  +0.98%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.55'
$ws.Range("E27").Value = '  +0.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.928'
$ws.Range("E28").Value = '  -0.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.434'
$ws.Range("E29").Value = '  -3.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09197'
$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("E31").Value = '  -2.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.964'
$ws.Range("E32").Value = '  -3.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04960'
$ws.Range("E33").Value = '  -3.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.143'
$ws.Range("E34").Value = '  +3.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7230'
$ws.Range("E35").Value = '  -2.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.685'
$ws.Range("E36").Value = '  -1.58%  '

$ws.Range("E37").Value = '  -0.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.596'
$ws.Range("E38").Value = '  -3.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8984'
$ws.Range("E39").Value = '  -2.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.903'
$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.99'
$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.48%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4239'
$ws.Range("E44").Value = '  -3.90%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.273'
$ws.Range("E45").Value = '  -4.00%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1299'
$ws.Range("E46").Value = '  -5.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '63.58'
$ws.Range("E47").Value = '  -8.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.490'
$ws.Range("E48").Value = '  +6.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.710'
$ws.Range("E49").Value = '  -4.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.77'
$ws.Range("E50").Value = '  -3.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05638'
$ws.Range("E51").Value = '  -3.39%  '
